$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values stay as text (avoid Excel numeric auto-conversion),
# matching the inline-string storage already used for this column.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "58.405.65"
$ws.Range("E2").Value = "  -3.94%  "

$ws.Range("D3").Value = "2.647.84"
$ws.Range("E3").Value = "  -2.24%  "

$ws.Range("E4").Value = "  +0.25%  "

$ws.Range("D5").Value = "521.72"
$ws.Range("E5").Value = "  -0.68%  "

$ws.Range("D6").Value = "144.85"
$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("E7").Value = "  +0.20%  "

$ws.Range("E8").Value = "  -1.12%  "

$ws.Range("E9").Value = "  +4.14%  "

$ws.Range("E10").Value = "  -3.09%  "

$ws.Range("D11").Value = "0.340"
$ws.Range("E11").Value = "  +0.42%  "

$ws.Range("E12").Value = "  +1.66%  "

$ws.Range("D13").Value = "3.118.24"
$ws.Range("E13").Value = "  -1.23%  "

$ws.Range("D14").Value = "58.417.19"
$ws.Range("E14").Value = "  -3.78%  "

$ws.Range("D15").Value = "20.89"
$ws.Range("E15").Value = "  -2.13%  "

$ws.Range("D16").Value = "0.0000136"
$ws.Range("E16").Value = "  -1.43%  "

$ws.Range("D17").Value = "2.656.32"
$ws.Range("E17").Value = "  -1.75%  "

$ws.Range("D18").Value = "338.87"
$ws.Range("E18").Value = "  -3.09%  "

$ws.Range("D19").Value = "4.39"
$ws.Range("E19").Value = "  -2.86%  "

$ws.Range("D20").Value = "10.43"
$ws.Range("E20").Value = "  -1.44%  "

$ws.Range("D21").Value = "6.33"
$ws.Range("E21").Value = "  +0.13%  "

$ws.Range("E22").Value = "  +0.07%  "

$ws.Range("D23").Value = "64.41"
$ws.Range("E23").Value = "  +1.11%  "

$ws.Range("E24").Value = "  +1.06%  "

$ws.Range("D25").Value = "0.166"
$ws.Range("E25").Value = "  -2.04%  "

$ws.Range("E26").Value = "  +0.46%  "

$ws.Range("D27").Value = "0.0₃0798"
$ws.Range("E27").Value = "  -2.72%  "

$ws.Range("D28").Value = "7.15"
$ws.Range("E28").Value = "  -2.60%  "

$ws.Range("D29").Value = "6.65"
$ws.Range("E29").Value = "  -3.83%  "

$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.13%  "

$ws.Range("D31").Value = "1.59"
$ws.Range("E31").Value = "  -0.57%  "

$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "152.00"
$ws.Range("E32").Value = "  +1.13%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "18.90"
$ws.Range("E33").Value = "  -2.02%  "

$ws.Range("D34").Value = "4.17"
$ws.Range("E34").Value = "  -2.16%  "

$ws.Range("D35").Value = "0.914"
$ws.Range("E35").Value = "  -4.05%  "

$ws.Range("E36").Value = "  -5.68%  "

$ws.Range("D37").Value = "0.870"
$ws.Range("E37").Value = "  -1.59%  "

$ws.Range("D38").Value = "36.84"
$ws.Range("E38").Value = "  -0.20%  "

$ws.Range("E39").Value = "  -5.11%  "

$ws.Range("D40").Value = "3.62"
$ws.Range("E40").Value = "  -1.02%  "

$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("D42").Value = "0.610"
$ws.Range("E42").Value = "  -0.21%  "

$ws.Range("D43").Value = "274.10"
$ws.Range("E43").Value = "  -3.23%  "

$ws.Range("E44").Value = "  -1.89%  "

$ws.Range("D45").Value = "19.46"
$ws.Range("E45").Value = "  -3.45%  "

$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").Value = "10.62"
$ws.Range("E46").Value = "  +1.62%  "

$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "0.0535"
$ws.Range("E47").Value = "  -0.88%  "

$ws.Range("D48").Value = "2.042.69"
$ws.Range("E48").Value = "  -4.70%  "

$ws.Range("E49").Value = "  -4.30%  "

$ws.Range("D50").Value = "0.0228"
$ws.Range("E50").Value = "  -2.89%  "

$ws.Range("D51").Value = "18.40"
$ws.Range("E51").Value = "  -3.06%  "
